$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each pair is a cell reference and its new text value. Coinranking prices,
# volumes and percentage changes are stored as literal text in this sheet
# (t="inlineStr"/shared-string cells, never numeric cells), and three rows
# (B7:E19) also rotate which coin/link occupies which row. We force every
# cell to Text format before writing so Excel does not auto-convert the
# numeric-looking strings (prices, "x.xx%") into real numbers, then clear
# the temporary formatting again so no stray cell style is left behind.
$updates = @(
    @("D2", "287.05"),
    @("E2", "1.32%"),
    @("E3", "3.86%"),
    @("E4", "1.50%"),
    @("D5", "0.06745"),
    @("E5", "3.76%"),
    @("D6", "7.348"),
    @("E6", "1.68%"),
    @("B7", "FTXToken"),
    @("C7", "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"),
    @("D7", "1.385"),
    @("E7", "-1.08%"),
    @("B8", "MXToken"),
    @("C8", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"),
    @("D8", "0.9173"),
    @("E8", "-0.06%"),
    @("B9", "WazirX"),
    @("C9", "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"),
    @("D9", "0.1592"),
    @("E9", "3.09%"),
    @("B10", "LiechtensteinCryptoassetsExchange"),
    @("C10", "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"),
    @("D10", "0.06846"),
    @("E10", "7.72%"),
    @("B11", "MandalaExchangeToken"),
    @("C11", "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"),
    @("D11", "0.07663"),
    @("E11", "1.32%"),
    @("B12", "BitrueCoin"),
    @("C12", "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"),
    @("D12", "0.02916"),
    @("E12", "2.09%"),
    @("B13", "BitMartToken"),
    @("C13", "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"),
    @("D13", "0.08983"),
    @("E13", "-0.09%"),
    @("B14", "BitForexToken"),
    @("C14", "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"),
    @("D14", "0.001584"),
    @("E14", "-0.63%"),
    @("B15", "CoinExToken"),
    @("C15", "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"),
    @("D15", "0.04468"),
    @("E15", "0.72%"),
    @("B16", "One"),
    @("C16", "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"),
    @("D16", "0.0006461"),
    @("E16", "1.30%"),
    @("B17", "TigerCash"),
    @("C17", "https://coinranking.com/coin/6hIn06L2+tigercash-tch"),
    @("D17", "0.006265"),
    @("E17", "1.51%"),
    @("B18", "LEO"),
    @("C18", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"),
    @("D18", "3.451"),
    @("E18", "0.26%"),
    @("B19", "GateToken"),
    @("C19", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"),
    @("D19", "3.447"),
    @("E19", "2.58%"),
    @("D21", "0.3200"),
    @("E21", "0.56%"),
    @("E22", "-1.17%"),
    @("D23", "4.060"),
    @("E23", "2.03%"),
    @("D25", "0.001191"),
    @("E25", "0.79%"),
    @("D26", "0.004138"),
    @("E26", "-7.04%"),
    @("D27", "0.0001196"),
    @("E27", "-0.33%"),
    @("E28", "-0.15%"),
    @("D40", "0.04265"),
    @("E40", "3.65%"),
    @("D41", "0.006727"),
    @("E41", "1.06%"),
    @("D42", "0.1243"),
    @("E42", "0.85%"),
    @("D43", "0.002142"),
    @("E43", "0.14%"),
    @("D44", "0.01196"),
    @("E44", "3.69%"),
    @("D45", "0.00005693"),
    @("E45", "1.63%"),
    @("D46", "1.963"),
    @("E46", "0.33%")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.ClearFormats()
}
